# "Added dynamic light support"
#
# Marks the dynamic-light rubric rows (directional / point / spot light,
# plus their "dynamic change" variants) as achieved at Milestone II.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows completed for Milestone II ("II" in the Student column, "X" in the
# Milestone Completed column). Row 28 (combining two lights) is left as-is.
$rows = @(25, 26, 27, 29, 30, 31)
foreach ($r in $rows) {
    $ws.Range("E$r").Value = "II"
    $ws.Range("F$r").Value = "X"
}

# Restore the on-screen scroll position / selection to match the author's
# saved view (scrolled down to row 10, with E18 selected).
$ws.Range("E18").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1

$wb.Application.Calculate()
